# Remove the personnel record for "AHMET YILDIRIM" (row with serial 26352,
# K.Teknisyeni, Silifke, phone 3242801542) from the MERSIN sheet.
# Deleting the entire row shifts all subsequent rows up by one and Excel
# automatically drops the now-unused "AHMET YILDIRIM" entry from the shared
# string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row to delete by its known value, rather than a hard-coded
# row number, to be robust.
$found = $ws.Columns.Item(1).Find("AHMET YILDIRIM")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

# The used range shrank from 99 to 98 rows; refresh the AutoFilter so its
# range reference (and the associated _FilterDatabase defined name) follow
# the new extent instead of keeping the stale E1:E99 range.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$ws.AutoFilterMode = $false
$ws.Range("E1:E$lastRow").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=MERSİN!`$E`$1:`$E`$$lastRow"
    }
}

# Match the saved cursor position left behind in the source file.
$ws.Range("I39").Select() | Out-Null
